$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value2 = 1569.5625
$ws.Range("I111").Value2 = 1820.2727
$ws.Range("J111").Value2 = 1018
$ws.Range("K111").Value2 = 5460.8181
$ws.Range("L111").Value2 = 3054
$ws.Range("M111").Value2 = -2393.8181
$ws.Range("N111").Value2 = -9188
$ws.Range("H132").Value2 = 2364.5833
$ws.Range("I132").Value2 = 2364.5833
$ws.Range("K132").Value2 = 7093.749899999999
$ws.Range("M132").Value2 = -4563.749899999999
$ws.Range("H137").Value2 = 2509.238
$ws.Range("I137").Value2 = 2835.3103
$ws.Range("J137").Value2 = 1781.8462
$ws.Range("K137").Value2 = 8505.930899999999
$ws.Range("L137").Value2 = 5345.5386
$ws.Range("M137").Value2 = -5955.930899999999
$ws.Range("N137").Value2 = -10445.5386
$ws.Range("H138").Value2 = 1969.98
$ws.Range("J138").Value2 = 2178.3901
$ws.Range("L138").Value2 = 6535.1703
$ws.Range("N138").Value2 = -16815.1703
$ws.Range("H141").Value2 = 10462.786
$ws.Range("I141").Value2 = 11410.833
$ws.Range("K141").Value2 = 34232.499
$ws.Range("M141").Value2 = -29052.499

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 13167960
$ws.Range("I32").Value2 = 16137466
$ws.Range("K32").Value2 = 16137466
$ws.Range("M32").Value2 = -16137179
$ws.Range("H35").Value2 = 4666.3335
$ws.Range("I35").Value2 = 4666.3335
$ws.Range("K35").Value2 = 4666.3335
$ws.Range("M35").Value2 = -4260.3335
$ws.Range("H45").Value2 = 38463740
$ws.Range("J45").Value2 = 4260.5
$ws.Range("L45").Value2 = 4260.5
$ws.Range("N45").Value2 = -5014.5
$ws.Range("H61").Value2 = 39478280
$ws.Range("I61").Value2 = 29416606
$ws.Range("K61").Value2 = 29416606
$ws.Range("M61").Value2 = -29416394
$ws.Range("H74").Value2 = 10834207
$ws.Range("I74").Value2 = 12500802
$ws.Range("K74").Value2 = 12500802
$ws.Range("M74").Value2 = -12499928
$ws.Range("H77").Value2 = 10834207
$ws.Range("I77").Value2 = 12500802
$ws.Range("K77").Value2 = 62504010
$ws.Range("M77").Value2 = -62499642
$ws.Range("H120").Value2 = 70000
$ws.Range("J120").Value2 = 70000
$ws.Range("L120").Value2 = 70000
$ws.Range("N120").Value2 = -79676
$ws.Range("H122").Value2 = 3543.3572
$ws.Range("I122").Value2 = 2586.625
$ws.Range("J122").Value2 = 4819
$ws.Range("K122").Value2 = 7759.875
$ws.Range("L122").Value2 = 14457
$ws.Range("M122").Value2 = -5309.875
$ws.Range("N122").Value2 = -19357
$ws.Range("H132").Value2 = 3921.25
$ws.Range("I132").Value2 = 2176.05
$ws.Range("J132").Value2 = 12647.25
$ws.Range("K132").Value2 = 6528.150000000001
$ws.Range("L132").Value2 = 37941.75
$ws.Range("M132").Value2 = -3998.150000000001
$ws.Range("N132").Value2 = -43001.75
$ws.Range("H136").Value2 = 39478280
$ws.Range("I136").Value2 = 29416606
$ws.Range("K136").Value2 = 88249818
$ws.Range("M136").Value2 = -88247268

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 5050.5454
$ws.Range("J20").Value2 = 4131.3335
$ws.Range("L20").Value2 = 4131.3335
$ws.Range("N20").Value2 = -4625.3335
$ws.Range("H134").Value2 = 502189.66
$ws.Range("I134").Value2 = 2209.9443
$ws.Range("K134").Value2 = 6629.8329
$ws.Range("M134").Value2 = -4094.8329

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 2275.1052
$ws.Range("I132").Value2 = 1601.8823
$ws.Range("K132").Value2 = 4805.6469
$ws.Range("M132").Value2 = -2275.6469
$ws.Range("H134").Value2 = 1617.3334
$ws.Range("I134").Value2 = 1554.2858
$ws.Range("K134").Value2 = 4662.857400000001
$ws.Range("M134").Value2 = -2127.857400000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value2 = 1525.6
$ws.Range("I46").Value2 = 971.4
$ws.Range("J46").Value2 = 2079.8
$ws.Range("K46").Value2 = 2914.2
$ws.Range("L46").Value2 = 6239.400000000001
$ws.Range("M46").Value2 = -2823.2
$ws.Range("N46").Value2 = -6421.400000000001
$ws.Range("H50").Value2 = 455.10526
$ws.Range("I50").Value2 = 349
$ws.Range("J50").Value2 = 461
$ws.Range("K50").Value2 = 1047
$ws.Range("L50").Value2 = 1383
$ws.Range("M50").Value2 = -566
$ws.Range("N50").Value2 = -2345
$ws.Range("H53").Value2 = 455.10526
$ws.Range("I53").Value2 = 349
$ws.Range("J53").Value2 = 461
$ws.Range("K53").Value2 = 1047
$ws.Range("L53").Value2 = 1383
$ws.Range("M53").Value2 = -566
$ws.Range("N53").Value2 = -2345
$ws.Range("H105").Value2 = 9250
$ws.Range("J105").Value2 = 9250
$ws.Range("L105").Value2 = 27750
$ws.Range("N105").Value2 = -32992
$ws.Range("H120").Value2 = 15547.637
$ws.Range("I120").Value2 = 7670.6665
$ws.Range("J120").Value2 = 25000
$ws.Range("K120").Value2 = 23011.9995
$ws.Range("L120").Value2 = 75000
$ws.Range("M120").Value2 = -18173.9995
$ws.Range("N120").Value2 = -84676

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value2 = 9
$ws.Range("J17").Value2 = 9
$ws.Range("L17").Value2 = 9
$ws.Range("N17").Value2 = -345
$ws.Range("H49").Value2 = 36494.5
$ws.Range("J49").Value2 = 36494.5
$ws.Range("L49").Value2 = 36494.5
$ws.Range("N49").Value2 = -36862.5
$ws.Range("H132").Value2 = 125016040
$ws.Range("I132").Value2 = 250001900
$ws.Range("J132").Value2 = 30170.75
$ws.Range("K132").Value2 = 750005700
$ws.Range("L132").Value2 = 90512.25
$ws.Range("M132").Value2 = -750003170
$ws.Range("N132").Value2 = -95572.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 51012.91
$ws.Range("I7").Value2 = 4164.75
$ws.Range("J7").Value2 = 107230.7
$ws.Range("K7").Value2 = 4164.75
$ws.Range("L7").Value2 = 107230.7
$ws.Range("M7").Value2 = -4052.75
$ws.Range("N7").Value2 = -107454.7
$ws.Range("H22").Value2 = 1370.3636
$ws.Range("I22").Value2 = 1275
$ws.Range("J22").Value2 = 1424.8572
$ws.Range("K22").Value2 = 1275
$ws.Range("L22").Value2 = 1424.8572
$ws.Range("M22").Value2 = -980
$ws.Range("N22").Value2 = -2014.8572
$ws.Range("H27").Value2 = 1370.3636
$ws.Range("I27").Value2 = 1275
$ws.Range("J27").Value2 = 1424.8572
$ws.Range("K27").Value2 = 1275
$ws.Range("L27").Value2 = 1424.8572
$ws.Range("M27").Value2 = -1168
$ws.Range("N27").Value2 = -1638.8572
$ws.Range("H122").Value2 = 4992.343
$ws.Range("I122").Value2 = 4681.227
$ws.Range("K122").Value2 = 14043.681
$ws.Range("M122").Value2 = -11593.681
$ws.Range("H126").Value2 = 51012.91
$ws.Range("I126").Value2 = 4164.75
$ws.Range("J126").Value2 = 107230.7
$ws.Range("K126").Value2 = 12494.25
$ws.Range("L126").Value2 = 321692.1
$ws.Range("M126").Value2 = -10024.25
$ws.Range("N126").Value2 = -326632.1

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value2 = 67995
$ws.Range("J127").Value2 = 67995
$ws.Range("L127").Value2 = 67995
$ws.Range("N127").Value2 = -77915
$ws.Range("H132").Value2 = 7028.2666
$ws.Range("I132").Value2 = 7775.857
$ws.Range("K132").Value2 = 23327.571
$ws.Range("M132").Value2 = -20797.571
$ws.Range("H141").Value2 = 64333.332
$ws.Range("J141").Value2 = 64333.332
$ws.Range("L141").Value2 = 64333.332
$ws.Range("N141").Value2 = -74693.33199999999
